$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B5 (TreatmentTab query): drop the redundant CONCAT(...) wrapper around REPLACE(...) ---
$b5 = $ws.Range("B5")
$oldText = $b5.Text
$newText = $oldText.Replace("CONCAT(REPLACE(trt.treatment_agent, ';', ', '))", "REPLACE(trt.treatment_agent, ';', ', ')")
$b5.Value = $newText

# Bump the font size on B5 (re-applied via the UI), which is what produced the
# extra font/style entries in the saved workbook.
$b5.Font.Size = 12

# --- Selection / view bookkeeping to match the saved workbook ---
$ws.Range("B2").Select()
